{"js": "// The commit splits the title run\n//   \"An\u00e1lisis Previos de Cambios para TP2\"\n// into three runs that all share the exact same run formatting as the\n// original run:\n//   1) \"An\u00e1lisis Previos de Cambios\"\n//   2) \" Previos\"\n//   3) \" para TP2\"\n// (the visible text becomes \"An\u00e1lisis Previos de Cambios Previos para TP2\").\n\nconst originalText = \"An\u00e1lisis Previos de Cambios para TP2\";\nconst firstPart = \"An\u00e1lisis Previos de Cambios\";\nconst secondPart = \" Previos\";\nconst thirdPart = \" para TP2\";\n\n// Find the range holding the original title text.\nconst body = context.document.body;\nconst results = body.search(originalText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the title text to split.\");\n}\nconst range = results.items[0];\n\n// Read back the run's own OOXML so we can reuse its *exact* run formatting\n// (rFonts/sz/szCs/u/lang \u2026) for the two brand-new runs we are about to add,\n// instead of hard-coding it.\nconst ooxmlResult = range.getOoxml();\nawait context.sync();\nconst rPrMatch = ooxmlResult.value.match(/<w:rPr>[\\s\\S]*?<\\/w:rPr>/);\nconst rPr = rPrMatch ? rPrMatch[0] : \"<w:rPr/>\";\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nfunction runXml(text, preserveSpace) {\n  const space = preserveSpace ? ' xml:space=\"preserve\"' : \"\";\n  return `<w:r>${rPr}<w:t${space}>${escapeXml(text)}</w:t></w:r>`;\n}\n\n// Three sibling runs, replacing the single original run, but leaving the\n// paragraph's own <w:pPr> (alignment, spacing, \u2026) untouched since we only\n// target the text range, not the whole paragraph.\nconst runsXml = runXml(firstPart, false) + runXml(secondPart, true) + runXml(thirdPart, true);\n\nconst packageXml =\n  `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>` +\n  `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n  `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n  `<pkg:xmlData>` +\n  `<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">` +\n  `<w:body><w:p>${runsXml}</w:p></w:body>` +\n  `</w:document>` +\n  `</pkg:xmlData></pkg:part></pkg:package>`;\n\nrange.insertOoxml(packageXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The commit splits the title run\n#   \"An\u00e1lisis Previos de Cambios para TP2\"\n# into three runs that all share the exact same run formatting as the\n# original run:\n#   1) \"An\u00e1lisis Previos de Cambios\"\n#   2) \" Previos\"\n#   3) \" para TP2\"\n# (the visible text becomes \"An\u00e1lisis Previos de Cambios Previos para TP2\").\n\n$d = $word.ActiveDocument\n\n$originalText = \"An\u00e1lisis Previos de Cambios para TP2\"\n$firstPart    = \"An\u00e1lisis Previos de Cambios\"\n$secondPart   = \" Previos\"\n$thirdPart    = \" para TP2\"\n\n# Locate the range holding the original title text.\n$targetRange = $d.Content\n$find = $targetRange.Find\n$find.Text = $originalText\n$find.MatchCase = $true\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find the title text to split.\"\n}\n\n# Read back the run's own OOXML so we can reuse its *exact* run formatting\n# (rFonts/sz/szCs/u/lang ...) for the two brand-new runs we are about to add,\n# instead of hard-coding it.\n$xml = $targetRange.WordOpenXML\nif ($xml -match '(?s)<w:rPr>.*?</w:rPr>') {\n    $rPr = $matches[0]\n} else {\n    $rPr = \"<w:rPr/>\"\n}\n\nfunction New-RunXml([string]$text, [bool]$preserve) {\n    if ($preserve) {\n        $space = ' xml:space=\"preserve\"'\n    } else {\n        $space = ''\n    }\n    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n    return \"<w:r>$rPr<w:t$space>$escaped</w:t></w:r>\"\n}\n\n$runsXml = (New-RunXml $firstPart $false) + (New-RunXml $secondPart $true) + (New-RunXml $thirdPart $true)\n\n# Replacing just this text range (not the whole paragraph) keeps the\n# paragraph's own <w:pPr> (alignment, spacing, ...) untouched.\n$packageXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>$runsXml</w:p></w:body>\" +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n$targetRange.InsertXML($packageXml)\n"}
